# Generate Report for Handback
# Updates the zh-cn and de-de report sheets with the new handback-file
# verification results for row 8 (8aa2d930-8369-4bd6-b24e-ba143e9769e5.md):
#   - populates "Latest Handback File" (I8) with a hyperlink to the file,
#   - populates "Latest Target File" / translated xlf name (J8),
#   - records the new handback datetime (K8),
#   - records an "Error Detail" (P8) because the handback commit used is stale,
#   - widens column P (Error Detail) so the message is readable.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e61cdd258de024513d909ab35b1dd0b02f5fff66/e2e/8aa2d930-8369-4bd6-b24e-ba143e9769e5.md"
$handbackDisplay = "8aa2d930-8369-4bd6-b24e-ba143e9769e5.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8864d423341d77f59ac765f2e0e2a325dd810af1/e2e/8aa2d930-8369-4bd6-b24e-ba143e9769e5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e61cdd258de024513d909ab35b1dd0b02f5fff66/e2e/8aa2d930-8369-4bd6-b24e-ba143e9769e5.md."

# Hyperlink font colour (Excel colour is BGR-packed): RGB(0x64,0x95,0xED)
$hyperlinkColor = 15570276

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsZh.Range("I8").Font.Underline = 2
$wsZh.Range("I8").Font.Color = $hyperlinkColor

$wsZh.Range("J8").Value = "8aa2d930-8369-4bd6-b24e-ba143e9769e5.f24bc13c1207b828270d60da08c4cf947edb4597.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-27 00:43:28"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsDe.Range("I8").Font.Underline = 2
$wsDe.Range("I8").Font.Color = $hyperlinkColor

$wsDe.Range("J8").Value = "8aa2d930-8369-4bd6-b24e-ba143e9769e5.f24bc13c1207b828270d60da08c4cf947edb4597.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-27 00:43:35"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.14
